# daily auto push: 2026-01-07 18:49 UTC
# Appends the next two time-slot rows ("2026/01/07" 23:00 and "2026/01/08" 02:00)
# to the daily log table on Sheet1, inserting them right after the last existing
# "2026/01/07" row (row 576) and pushing the remaining rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 576:577 - everything currently at/after row 576
# (starting with the 2026/12/29 block) shifts down to 578 onward, matching the
# rest of the table's existing formatting (no explicit style).
$ws.Rows("576:577").Insert()

# Column A holds dates stored as plain text (e.g. "2026/12/29"), not real Excel
# dates. Assigning a literal "YYYY/MM/DD" string to .Value gets auto-parsed into
# a date value, and forcing text via NumberFormat="@" (or a leading apostrophe)
# leaves a quote-prefix style behind that the source cells don't have. Instead,
# write each date as a short text formula and flatten it to a plain value via
# copy / paste-special so the cell ends up as ordinary text with the default
# (unstyled) format, same as every other cell in the column.
$ws.Range("A576").Formula = "=""2026/01/07"""
$ws.Range("A577").Formula = "=""2026/01/08"""
$ws.Range("A576:A577").Copy()
$ws.Range("A576:A577").PasteSpecial(-4163)

# Row 576: 2026/01/07 (Wed), time slot 23, ranking 24
$ws.Range("B576").Value = "水"
$ws.Range("C576").Value = 23
$ws.Range("D576").Value = 24

# Row 577: 2026/01/08 (Thu), time slot 2, ranking 25
$ws.Range("B577").Value = "木"
$ws.Range("C577").Value = 2
$ws.Range("D577").Value = 25
